# Predetermination Letter revert edit
# Reverts an NPPF-2021 "keystone" update back to the NPPF-2019 wording,
# removes the stray _GoBack bookmark around "received on ", and fixes
# "Archaeology Adviser" -> "Archaeology Advisor".

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Title year: 2021 -> 2019
Replace-Text "NATIONAL PLANNING POLICY FRAMEWORK 2021" "NATIONAL PLANNING POLICY FRAMEWORK 2019"

# 2. Remove the "_GoBack" bookmark by merging the two runs it straddles
#    ("Thank you for your consultation " + "received on ") back into one
#    contiguous run of text.
Replace-Text "Thank you for your consultation received on " "Thank you for your consultation received on "

# 3. London Plan year: 2021 -> 2017
Replace-Text "the London Plan (2021 Policy HC1)" "the London Plan (2017 Policy HC1)"

# 4. NPPF paragraph numbers (2021 NPPF -> 2019 NPPF numbering)
Replace-Text "NPPF paragraph 194 says" "NPPF paragraph 189 says"
Replace-Text "NPPF paragraphs 199 - 202 place" "NPPF paragraphs 193 - 194 place"
Replace-Text "(NPPF paragraph 203)" "(NPPF paragraph 197)"
Replace-Text "NPPF paragraphs 190 and 197 and London Plan" "NPPF paragraphs 185 and 192 and London Plan"
Replace-Text "paragraph 205 of the NPPF" "paragraph 199 of the NPPF"

# 5. Job title typo fix
Replace-Text "Archaeology Adviser" "Archaeology Advisor"
